$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 96 (previously rank 94 -> now 98)
$ws.Range("A96").Value = 98
$ws.Range("B96").Value = 1.901710291787398
$ws.Range("C96").Value = 4185.8
$ws.Range("D96").Value = 0.01349240180076153
$ws.Range("E96").Value = 29.6
$ws.Range("F96").Value = 204.4
$ws.Range("G96").Value = "MyDogaN"
$ws.Range("H96").Value = "DUO_SUPPORT"
$ws.Range("I96").Value = 0.09128055152368671
$ws.Range("J96").Value = 18.6
$ws.Range("K96").Value = 0.008167217339014521

# Update row 97 (previously rank 96 -> now 99)
$ws.Range("A97").Value = 99
$ws.Range("B97").Value = 2.791646791513082
$ws.Range("C97").Value = 3890.6
$ws.Range("D97").Value = 0.01949175076877424
$ws.Range("E97").Value = 28.2
$ws.Range("F97").Value = 221.4
$ws.Range("G97").Value = "Mr Kayn"
$ws.Range("H97").Value = "DUO_SUPPORT"
$ws.Range("I97").Value = 0.1749684198889241
$ws.Range("J97").Value = 3.2
$ws.Range("K97").Value = 0.002463335941977546

# Update row 98 (previously rank 97 -> now 100)
$ws.Range("A98").Value = 100
$ws.Range("B98").Value = 2.416763848396501
$ws.Range("C98").Value = 3315.8
$ws.Range("D98").Value = 0.0163265306122449
$ws.Range("E98").Value = 22.4
$ws.Range("F98").Value = 40.6
$ws.Range("G98").Value = "Booogeyman"
$ws.Range("H98").Value = "DUO_CARRY"
$ws.Range("I98").Value = 0.02959183673469388
$ws.Range("J98").Value = 3.4
$ws.Range("K98").Value = 0.002478134110787172

# Update row 99 (previously rank 98 -> now 101)
$ws.Range("A99").Value = 101
$ws.Range("B99").Value = 9.48417331160671
$ws.Range("C99").Value = 15096.8
$ws.Range("D99").Value = 0.04245232838056302
$ws.Range("E99").Value = 67
$ws.Range("F99").Value = 161.6
$ws.Range("G99").Value = "Shiller"
$ws.Range("H99").Value = "DUO_CARRY"
$ws.Range("I99").Value = 0.1032920977440951
$ws.Range("J99").Value = 5.4
$ws.Range("K99").Value = 0.003379839878639451

# New row 100 (rank 102) - copy formatting of column A from the row above first
$ws.Range("A99").Copy()
$ws.Range("A100").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A100").Value = 102
$ws.Range("B100").Value = 1.249264705882353
$ws.Range("C100").Value = 2038.8
$ws.Range("D100").Value = 0.01409313725490196
$ws.Range("E100").Value = 23
$ws.Range("F100").Value = 138
$ws.Range("G100").Value = "Poppy Gods"
$ws.Range("H100").Value = "DUO_CARRY"
$ws.Range("I100").Value = 0.08455882352941177
$ws.Range("J100").Value = 3
$ws.Range("K100").Value = 0.001838235294117647

Write-Output "Edit applied"
